$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = 'UserID'
$ws.Range("B1").Value = 'PutwallPickingQuantity'
$ws.Range("C1").Value = 'UPH'
$ws.Range("A2").Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Range("B2").Value = 26
$ws.Range("C2").Value = 14.18181818181818
$ws.Range("A3").Value = 'BOHD0676.KUSHLIAK'
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1.090909090909091
$ws.Range("A4").Value = 'DIAN4065.ENTRIALGO'
$ws.Range("B4").Value = 67
$ws.Range("C4").Value = 36.54545454545455
$ws.Range("A5").Value = 'PATR5027.AMEH'
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1.636363636363636
$ws.Range("A6").Value = 'THIE6554.DIALLO'
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 7.636363636363637
$ws.Range("A7").Value = 'XUAN754N.LU'
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 1.636363636363636
$ws.Range("A8").Value = 'ZAHIDGUL.MINHAS'
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 2.181818181818182

$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = 'UserID'
$ws.Range("B1").Value = 'RegularPickQuantity'
$ws.Range("C1").Value = 'UPH'
$ws.Range("A2").Value = 'BOHD0676.KUSHLIAK'
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 2.181818181818182
$ws.Range("A3").Value = 'DIAN4065.ENTRIALGO'
$ws.Range("B3").Value = 49
$ws.Range("C3").Value = 26.72727272727273
$ws.Range("A4").Value = 'PATR5027.AMEH'
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1.090909090909091
$ws.Range("A5").Value = 'ZAHIDGUL.MINHAS'
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1.636363636363636

$ws = $wb.Worksheets.Item(4)
$ws.Range("A4:C7").Clear() | Out-Null
$ws.Range("A1").Value = 'UserID'
$ws.Range("B1").Value = 'SinglePickQuantity'
$ws.Range("C1").Value = 'UPH'
$ws.Range("A2").Value = 'REJWAN.ISLAM'
$ws.Range("B2").Value = 18
$ws.Range("C2").Value = 9.818181818181818
$ws.Range("A3").Value = 'WESL5337.CADETTE'
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 22.36363636363636

$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = 'UserID'
$ws.Range("B1").Value = 'ReplenishmentPickQuantity'
$ws.Range("C1").Value = 'UPH'
$ws.Range("A2").Value = 'AHME0710.JUBRAN'
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 5.454545454545455
$ws.Range("A3").Value = 'ANASTASIIA.MAKHTOUT'
$ws.Range("B3").Value = 93
$ws.Range("C3").Value = 50.72727272727273
$ws.Range("A4").Value = 'ANJALI.BAKSHI'
$ws.Range("B4").Value = 94
$ws.Range("C4").Value = 51.27272727272727
$ws.Range("A5").Value = 'BUDD0680.TENNAKOON'
$ws.Range("B5").Value = 79
$ws.Range("C5").Value = 43.09090909090909
$ws.Range("A6").Value = 'DEVI789.SINGH'
$ws.Range("B6").Value = 35
$ws.Range("C6").Value = 19.09090909090909
$ws.Range("A7").Value = 'GIGNESH.PATEL'
$ws.Range("B7").Value = 34
$ws.Range("C7").Value = 18.54545454545455
$ws.Range("A8").Value = 'IREN797N.CABRERA'
$ws.Range("B8").Value = 65
$ws.Range("C8").Value = 35.45454545454545
$ws.Range("A9").Value = 'JEEW9554.SITUMUDALIG'
$ws.Range("B9").Value = 132
$ws.Range("C9").Value = 72
$ws.Range("A10").Value = 'KADE3054.ZONGO'
$ws.Range("B10").Value = 51
$ws.Range("C10").Value = 27.81818181818182
$ws.Range("A11").Value = 'LOWRHY-OTIENO.JAOKO'
$ws.Range("B11").Value = 99
$ws.Range("C11").Value = 54
$ws.Range("A12").Value = 'MICA0432.RIZKALLAMAR'
$ws.Range("B12").Value = 111
$ws.Range("C12").Value = 60.54545454545455
$ws.Range("A13").Value = 'PATI2298.ATSIANGBE'
$ws.Range("B13").Value = 24
$ws.Range("C13").Value = 13.09090909090909
$ws.Range("A14").Value = 'PRINCE.FORSON'
$ws.Range("B14").Value = 48
$ws.Range("C14").Value = 26.18181818181818
$ws.Range("A15").Value = 'REJWAN.ISLAM'
$ws.Range("B15").Value = 40
$ws.Range("C15").Value = 21.81818181818182
$ws.Range("A16").Value = 'STAN9294.BAUER'
$ws.Range("B16").Value = 61
$ws.Range("C16").Value = 33.27272727272727
$ws.Range("A17").Value = 'SURESH.DHAWAN'
$ws.Range("B17").Value = 36
$ws.Range("C17").Value = 19.63636363636364
$ws.Range("A18").Value = 'THIE6554.DIALLO'
$ws.Range("B18").Value = 70
$ws.Range("C18").Value = 38.18181818181818
$ws.Range("A19").Value = 'WESL5337.CADETTE'
$ws.Range("B19").Value = 62
$ws.Range("C19").Value = 33.81818181818182
$ws.Range("A20").Value = 'XUAN754N.LU'
$ws.Range("B20").Value = 59
$ws.Range("C20").Value = 32.18181818181818
$ws.Range("A21").Value = 'YATI0689.YATIN'
$ws.Range("B21").Value = 38
$ws.Range("C21").Value = 20.72727272727273
$ws.Range("A22").Value = 'ZAKI0190.PHILLIPHORS'
$ws.Range("B22").Value = 115
$ws.Range("C22").Value = 62.72727272727273

$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = 'UserID'
$ws.Range("B1").Value = 'QuickMoveQuantity'
$ws.Range("C1").Value = 'UPH'
$ws.Range("A2").Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Range("B2").Value = 115
$ws.Range("C2").Value = 62.72727272727273
$ws.Range("A3").Value = 'BOHD0676.KUSHLIAK'
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1.090909090909091
$ws.Range("A4").Value = 'DEVI789.SINGH'
$ws.Range("B4").Value = 27
$ws.Range("C4").Value = 14.72727272727273
$ws.Range("A5").Value = 'DIAN4065.ENTRIALGO'
$ws.Range("B5").Value = 59
$ws.Range("C5").Value = 32.18181818181818
$ws.Range("A6").Value = 'ESSE0616.UDEH'
$ws.Range("B6").Value = 65
$ws.Range("C6").Value = 35.45454545454545
$ws.Range("A7").Value = 'KADE3054.ZONGO'
$ws.Range("B7").Value = 28
$ws.Range("C7").Value = 15.27272727272727
$ws.Range("A8").Value = 'MICA0432.RIZKALLAMAR'
$ws.Range("B8").Value = 111
$ws.Range("C8").Value = 60.54545454545455
$ws.Range("A9").Value = 'NESR2403.ATTALAH'
$ws.Range("B9").Value = 31
$ws.Range("C9").Value = 16.90909090909091
$ws.Range("A10").Value = 'STAN9294.BAUER'
$ws.Range("B10").Value = 60
$ws.Range("C10").Value = 32.72727272727273
$ws.Range("A11").Value = 'SURESH.DHAWAN'
$ws.Range("B11").Value = 123
$ws.Range("C11").Value = 67.09090909090909
$ws.Range("A12").Value = 'THIE6554.DIALLO'
$ws.Range("B12").Value = 159
$ws.Range("C12").Value = 86.72727272727273
$ws.Range("A13").Value = 'WESL5337.CADETTE'
$ws.Range("B13").Value = 61
$ws.Range("C13").Value = 33.27272727272727
$ws.Range("A14").Value = 'XUAN754N.LU'
$ws.Range("B14").Value = 59
$ws.Range("C14").Value = 32.18181818181818
$ws.Range("A15").Value = 'YATI0689.YATIN'
$ws.Range("B15").Value = 97
$ws.Range("C15").Value = 52.90909090909091

$ws = $wb.Worksheets.Item(7)
$ws.Range("A27:B27").Clear() | Out-Null
$ws.Range("A1").Value = 'UserID'
$ws.Range("B1").Value = 'TotalIdleTime'
$ws.Range("A2").Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Range("B2").Value = 78
$ws.Range("A3").Value = 'AHME0710.JUBRAN'
$ws.Range("B3").Value = 49
$ws.Range("A4").Value = 'ANASTASIIA.MAKHTOUT'
$ws.Range("B4").Value = 50
$ws.Range("A5").Value = 'ANJALI.BAKSHI'
$ws.Range("B5").Value = 42
$ws.Range("A6").Value = 'BOHD0676.KUSHLIAK'
$ws.Range("B6").Value = 90
$ws.Range("A7").Value = 'BUDD0680.TENNAKOON'
$ws.Range("B7").Value = 27
$ws.Range("A8").Value = 'DEVI789.SINGH'
$ws.Range("B8").Value = 44
$ws.Range("A9").Value = 'DIAN4065.ENTRIALGO'
$ws.Range("B9").Value = 37
$ws.Range("A10").Value = 'ESSE0616.UDEH'
$ws.Range("B10").Value = 93
$ws.Range("A11").Value = 'GIGNESH.PATEL'
$ws.Range("B11").Value = 99
$ws.Range("A12").Value = 'IREN797N.CABRERA'
$ws.Range("B12").Value = 53
$ws.Range("A13").Value = 'JEEW9554.SITUMUDALIG'
$ws.Range("B13").Value = 21
$ws.Range("A14").Value = 'LOWRHY-OTIENO.JAOKO'
$ws.Range("B14").Value = 27
$ws.Range("A15").Value = 'NESR2403.ATTALAH'
$ws.Range("B15").Value = 109
$ws.Range("A16").Value = 'PATI2298.ATSIANGBE'
$ws.Range("B16").Value = 39
$ws.Range("A17").Value = 'PATR5027.AMEH'
$ws.Range("B17").Value = 64
$ws.Range("A18").Value = 'PRINCE.FORSON'
$ws.Range("B18").Value = 67
$ws.Range("A19").Value = 'REJWAN.ISLAM'
$ws.Range("B19").Value = 70
$ws.Range("A20").Value = 'SURESH.DHAWAN'
$ws.Range("B20").Value = 57
$ws.Range("A21").Value = 'THIE6554.DIALLO'
$ws.Range("B21").Value = 18
$ws.Range("A22").Value = 'WESL5337.CADETTE'
$ws.Range("B22").Value = 17
$ws.Range("A23").Value = 'XUAN754N.LU'
$ws.Range("B23").Value = 21
$ws.Range("A24").Value = 'YATI0689.YATIN'
$ws.Range("B24").Value = 43
$ws.Range("A25").Value = 'ZAHIDGUL.MINHAS'
$ws.Range("B25").Value = 94
$ws.Range("A26").Value = 'ZAKI0190.PHILLIPHORS'
$ws.Range("B26").Value = 28

$ws = $wb.Worksheets.Item(8)
$ws.Range("A1").Value = 'Hour'
$ws.Range("B1").Value = 'Regular Pick'
$ws.Range("C1").Value = 'Single Pick'
$ws.Range("D1").Value = 'Replenishment Pick'
$ws.Range("E1").Value = 'Putwall Pick'
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = -26
$ws.Range("C2").Value = -32
$ws.Range("D2").Value = -594
$ws.Range("E2").Value = -1
$ws.Range("A3").Value = 21
$ws.Range("B3").Value = -32
$ws.Range("C3").Value = -24
$ws.Range("D3").Value = -751
$ws.Range("E3").Value = -106
$ws.Range("A4").Value = 22
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = -3
$ws.Range("D4").Value = -11
$ws.Range("E4").Value = -12
$ws.Range("A5").Value = 23
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("A6").Value = 'Total'
$ws.Range("B6").Value = -58
$ws.Range("C6").Value = -59
$ws.Range("D6").Value = -1356
$ws.Range("E6").Value = -119
